$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert a new column before AV (shifts AV.. right by 1)
$ws.Columns("AV").Insert()

# Step 2: insert 6 new columns before the (now-shifted) BG, i.e. before old BF
$ws.Range("BG1:BL26").Insert(-4161)

# Step 3: header row - copy header style (bold/border/center) from B1 into new header cells
$ws.Range("B1").Copy()
$ws.Range("AV1").PasteSpecial(-4122)
$ws.Range("BG1").PasteSpecial(-4122)
$ws.Range("BH1").PasteSpecial(-4122)
$ws.Range("BI1").PasteSpecial(-4122)
$ws.Range("BJ1").PasteSpecial(-4122)
$ws.Range("BK1").PasteSpecial(-4122)
$ws.Range("BL1").PasteSpecial(-4122)

# Step 4: header values
$ws.Range("AV1").Value = "150_2_0.10_0.001"
$ws.Range("BG1").Value = "200_4_0.10_0.001"
$ws.Range("BH1").Value = "200_4_0.10_0.005"
$ws.Range("BI1").Value = "200_4_0.10_0.010"
$ws.Range("BJ1").Value = "200_4_0.10_0.020"
$ws.Range("BK1").Value = "200_4_0.10_0.050"
$ws.Range("BL1").Value = "200_4_0.10_0.100"

# Step 5: data values for rows 2..26
$ws.Range("AV2").Value = 0.9997333216416557
$ws.Range("BG2").Value = 0.9969165914032172
$ws.Range("BH2").Value = 0.9441361237575534
$ws.Range("BI2").Value = 0.6615499477729663
$ws.Range("BJ2").Value = 0.1615535829676643
$ws.Range("BK2").Value = 0.03487929122281073
$ws.Range("BL2").Value = 0.01428342100979789
$ws.Range("AV3").Value = 0.9997333216416557
$ws.Range("BG3").Value = 0.9863886309697759
$ws.Range("BH3").Value = 0.9306586356799963
$ws.Range("BI3").Value = 0.8330394420467792
$ws.Range("BJ3").Value = 0.4874808429108185
$ws.Range("BK3").Value = 0.058006083766369
$ws.Range("BL3").Value = 0.01419496698673787
$ws.Range("AV4").Value = 0.9997333216416557
$ws.Range("BG4").Value = 0.9731214201969496
$ws.Range("BH4").Value = 0.9102090202406775
$ws.Range("BI4").Value = 0.7808354281723168
$ws.Range("BJ4").Value = 0.3611228911221068
$ws.Range("BK4").Value = 0.04585470136190287
$ws.Range("BL4").Value = 0.01145188990612546
$ws.Range("AV5").Value = 0.9994666432833114
$ws.Range("BG5").Value = 0.994646752344113
$ws.Range("BH5").Value = 0.9625875713947006
$ws.Range("BI5").Value = 0.9085258472653507
$ws.Range("BJ5").Value = 0.6804719299471677
$ws.Range("BK5").Value = 0.0836596804541747
$ws.Range("BL5").Value = 0.0167874856413665
$ws.Range("AV6").Value = 1
$ws.Range("BG6").Value = 0.9965280742354939
$ws.Range("BH6").Value = 0.94363459591424
$ws.Range("BI6").Value = 0.6557856909304134
$ws.Range("BJ6").Value = 0.1867712167802713
$ws.Range("BK6").Value = 0.09227686753790049
$ws.Range("BL6").Value = 0.04372167221789883
$ws.Range("AV7").Value = 1
$ws.Range("BG7").Value = 0.9969259282465214
$ws.Range("BH7").Value = 0.96403549645121
$ws.Range("BI7").Value = 0.8752216713866423
$ws.Range("BJ7").Value = 0.4434134851995332
$ws.Range("BK7").Value = 0.1773603484864924
$ws.Range("BL7").Value = 0.1211106271230779
$ws.Range("AV8").Value = 1
$ws.Range("BG8").Value = 0.99598266580908
$ws.Range("BH8").Value = 0.9653112249206492
$ws.Range("BI8").Value = 0.8945068258727146
$ws.Range("BJ8").Value = 0.6383223698317224
$ws.Range("BK8").Value = 0.07106611390483203
$ws.Range("BL8").Value = 0.01874757386993528
$ws.Range("AV9").Value = 0.9997333216416557
$ws.Range("BG9").Value = 0.9234016706646517
$ws.Range("BH9").Value = 0.7463103743356317
$ws.Range("BI9").Value = 0.4195094201676254
$ws.Range("BJ9").Value = 0.1865276659785612
$ws.Range("BK9").Value = 0.04842614711368523
$ws.Range("BL9").Value = 0.01357219024108321
$ws.Range("AV10").Value = 0.9997333216416557
$ws.Range("BG10").Value = 0.995453296308074
$ws.Range("BH10").Value = 0.930664046932816
$ws.Range("BI10").Value = 0.5655554148906075
$ws.Range("BJ10").Value = 0.1635960017673818
$ws.Range("BK10").Value = 0.1142123123879254
$ws.Range("BL10").Value = 0.0451716818290927
$ws.Range("AV11").Value = 1
$ws.Range("BG11").Value = 0.9959826470295405
$ws.Range("BH11").Value = 0.9655502890729598
$ws.Range("BI11").Value = 0.8938965923276447
$ws.Range("BJ11").Value = 0.6188263894715703
$ws.Range("BK11").Value = 0.06025379217139843
$ws.Range("BL11").Value = 0.0150646994717561
$ws.Range("AV12").Value = 1
$ws.Range("BG12").Value = 0.9958594791021533
$ws.Range("BH12").Value = 0.9574916720813305
$ws.Range("BI12").Value = 0.7165841909956773
$ws.Range("BJ12").Value = 0.2379056042389926
$ws.Range("BK12").Value = 0.05239781044722126
$ws.Range("BL12").Value = 0.01937221272534833
$ws.Range("AV13").Value = 0.9994666432833114
$ws.Range("BG13").Value = 0.9955846552207513
$ws.Range("BH13").Value = 0.9606392134056354
$ws.Range("BI13").Value = 0.9034927436479241
$ws.Range("BJ13").Value = 0.6902401240804785
$ws.Range("BK13").Value = 0.08683588040474816
$ws.Range("BL13").Value = 0.01760479673076406
$ws.Range("AV14").Value = 0.9994666432833114
$ws.Range("BG14").Value = 0.995452088795534
$ws.Range("BH14").Value = 0.9612352058413096
$ws.Range("BI14").Value = 0.8985169263372967
$ws.Range("BJ14").Value = 0.6856817286886158
$ws.Range("BK14").Value = 0.119020324869413
$ws.Range("BL14").Value = 0.03731895820239367
$ws.Range("AV15").Value = 0.9994666432833114
$ws.Range("BG15").Value = 0.9959880791077411
$ws.Range("BH15").Value = 0.9685238385463478
$ws.Range("BI15").Value = 0.9072535032001448
$ws.Range("BJ15").Value = 0.6803528086047463
$ws.Range("BK15").Value = 0.1105864339642547
$ws.Range("BL15").Value = 0.03572801447038757
$ws.Range("AV16").Value = 0.9994666432833114
$ws.Range("BG16").Value = 0.9950569671055998
$ws.Range("BH16").Value = 0.9657608771023063
$ws.Range("BI16").Value = 0.8971301958407158
$ws.Range("BJ16").Value = 0.6440467949846129
$ws.Range("BK16").Value = 0.07923464908684433
$ws.Range("BL16").Value = 0.01812368752076735
$ws.Range("AV17").Value = 0.9994666432833114
$ws.Range("BG17").Value = 0.9943774133031127
$ws.Range("BH17").Value = 0.9620376917075388
$ws.Range("BI17").Value = 0.9073370992174812
$ws.Range("BJ17").Value = 0.6725817311980025
$ws.Range("BK17").Value = 0.07711567251692911
$ws.Range("BL17").Value = 0.01472129229869013
$ws.Range("AV18").Value = 0.9994666432833114
$ws.Range("BG18").Value = 0.9914693431475509
$ws.Range("BH18").Value = 0.9464334970444814
$ws.Range("BI18").Value = 0.8857781737701302
$ws.Range("BJ18").Value = 0.6830425692026182
$ws.Range("BK18").Value = 0.08696728074618967
$ws.Range("BL18").Value = 0.01868898448916305
$ws.Range("AV19").Value = 0.9994666432833114
$ws.Range("BG19").Value = 0.991603354348142
$ws.Range("BH19").Value = 0.9481955743298468
$ws.Range("BI19").Value = 0.8880501084940781
$ws.Range("BJ19").Value = 0.6760232287959481
$ws.Range("BK19").Value = 0.08124682403607722
$ws.Range("BL19").Value = 0.02152530674122039
$ws.Range("AV20").Value = 0.9994666432833114
$ws.Range("BG20").Value = 0.9945127411435221
$ws.Range("BH20").Value = 0.9633567136117311
$ws.Range("BI20").Value = 0.9085258134377218
$ws.Range("BJ20").Value = 0.6744159969450421
$ws.Range("BK20").Value = 0.115974631917495
$ws.Range("BL20").Value = 0.05531748575464013
$ws.Range("AV21").Value = 0.9994666432833114
$ws.Range("BG21").Value = 0.9951841874316938
$ws.Range("BH21").Value = 0.9600857355652052
$ws.Range("BI21").Value = 0.8995182755634746
$ws.Range("BJ21").Value = 0.6877154090945778
$ws.Range("BK21").Value = 0.1585444262462216
$ws.Range("BL21").Value = 0.1137845180768132
$ws.Range("AV22").Value = 0.9994666432833114
$ws.Range("BG22").Value = 0.9943801202281477
$ws.Range("BH22").Value = 0.9589075190577238
$ws.Range("BI22").Value = 0.8990653192286614
$ws.Range("BJ22").Value = 0.6896015694245827
$ws.Range("BK22").Value = 0.08525657146673733
$ws.Range("BL22").Value = 0.01904636420008394
$ws.Range("AV23").Value = 0.9994666432833114
$ws.Range("BG23").Value = 0.9914706967960417
$ws.Range("BH23").Value = 0.9514731711330778
$ws.Range("BI23").Value = 0.8960096133954376
$ws.Range("BJ23").Value = 0.6808596517141402
$ws.Range("BK23").Value = 0.07941978801571954
$ws.Range("BL23").Value = 0.01848243717579944
$ws.Range("AV24").Value = 0.9994666432833114
$ws.Range("BG24").Value = 0.9945127411435221
$ws.Range("BH24").Value = 0.9615579318560857
$ws.Range("BI24").Value = 0.9086157971217571
$ws.Range("BJ24").Value = 0.6814793404326468
$ws.Range("BK24").Value = 0.09659172722922246
$ws.Range("BL24").Value = 0.03068153185591176
$ws.Range("AV25").Value = 0.9994666432833114
$ws.Range("BG25").Value = 0.9943801202281477
$ws.Range("BH25").Value = 0.9589127505181126
$ws.Range("BI25").Value = 0.8995794371265589
$ws.Range("BJ25").Value = 0.689958454636613
$ws.Range("BK25").Value = 0.08549453420890687
$ws.Range("BL25").Value = 0.01799030689602236
$ws.Range("AV26").Value = 0.9994666432833114
$ws.Range("BG26").Value = 0.9946467523441129
$ws.Range("BH26").Value = 0.962028262368767
$ws.Range("BI26").Value = 0.9062907553562378
$ws.Range("BJ26").Value = 0.6803365772971084
$ws.Range("BK26").Value = 0.08124438061865892
$ws.Range("BL26").Value = 0.01899125654136917
